# Auto-generated edit script: updates crafting-profit market data cells
# per the scheduled runner refresh (commit: "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 491.35715
$ws.Range("I18").Value = 239.91667
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 239.91667
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 44.08332999999999
$ws.Range("N18").Value = -2568
$ws.Range("H38").Value = 1946.0834
$ws.Range("I38").Value = 95.90000000000001
$ws.Range("J38").Value = 3267.6428
$ws.Range("K38").Value = 287.7
$ws.Range("L38").Value = 9802.928400000001
$ws.Range("M38").Value = 84.29999999999995
$ws.Range("N38").Value = -10546.9284
$ws.Range("H58").Value = 1315.3334
$ws.Range("I58").Value = 484.72726
$ws.Range("J58").Value = 1796.2106
$ws.Range("K58").Value = 1454.18178
$ws.Range("L58").Value = 5388.6318
$ws.Range("M58").Value = -1304.18178
$ws.Range("N58").Value = -5688.6318
$ws.Range("H70").Value = 1593
$ws.Range("I70").Value = 1654
$ws.Range("J70").Value = 1501.5
$ws.Range("K70").Value = 4962
$ws.Range("L70").Value = 4504.5
$ws.Range("M70").Value = -4692
$ws.Range("N70").Value = -5044.5
$ws.Range("H73").Value = 1593
$ws.Range("I73").Value = 1654
$ws.Range("J73").Value = 1501.5
$ws.Range("K73").Value = 4962
$ws.Range("L73").Value = 4504.5
$ws.Range("M73").Value = -4026
$ws.Range("N73").Value = -6376.5
$ws.Range("H107").Value = 2633
$ws.Range("I107").Value = 1697.25
$ws.Range("J107").Value = 4771.857
$ws.Range("K107").Value = 1697.25
$ws.Range("L107").Value = 4771.857
$ws.Range("M107").Value = 222.75
$ws.Range("N107").Value = -8611.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8641.092000000001
$ws.Range("I32").Value = 6757.2534
$ws.Range("K32").Value = 6757.2534
$ws.Range("M32").Value = -6470.2534
$ws.Range("H45").Value = 1298.4166
$ws.Range("I45").Value = 1286.375
$ws.Range("J45").Value = 1322.5
$ws.Range("K45").Value = 1286.375
$ws.Range("L45").Value = 1322.5
$ws.Range("M45").Value = -909.375
$ws.Range("N45").Value = -2076.5
$ws.Range("H61").Value = 45456056
$ws.Range("J61").Value = 3028.8
$ws.Range("L61").Value = 3028.8
$ws.Range("N61").Value = -3452.8
$ws.Range("H74").Value = 3026.5715
$ws.Range("I74").Value = 2827.2856
$ws.Range("J74").Value = 3126.2144
$ws.Range("K74").Value = 2827.2856
$ws.Range("L74").Value = 3126.2144
$ws.Range("M74").Value = -1953.2856
$ws.Range("N74").Value = -4874.2144
$ws.Range("H77").Value = 3026.5715
$ws.Range("I77").Value = 2827.2856
$ws.Range("J77").Value = 3126.2144
$ws.Range("K77").Value = 14136.428
$ws.Range("L77").Value = 15631.072
$ws.Range("M77").Value = -9768.428
$ws.Range("N77").Value = -24367.072
$ws.Range("H132").Value = 2503.532
$ws.Range("I132").Value = 1779.3448
$ws.Range("K132").Value = 5338.0344
$ws.Range("M132").Value = -2808.0344
$ws.Range("H136").Value = 45456056
$ws.Range("J136").Value = 3028.8
$ws.Range("L136").Value = 9086.400000000001
$ws.Range("N136").Value = -14186.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 931.5517
$ws.Range("I107").Value = 806.05884
$ws.Range("K107").Value = 806.05884
$ws.Range("M107").Value = 1113.94116

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125001790
$ws.Range("I16").Value = 142858910
$ws.Range("K16").Value = 142858910
$ws.Range("M16").Value = -142858623
$ws.Range("H31").Value = 1156.0845
$ws.Range("I31").Value = 1106.5082
$ws.Range("J31").Value = 1458.5
$ws.Range("K31").Value = 1106.5082
$ws.Range("L31").Value = 1458.5
$ws.Range("M31").Value = -811.5082
$ws.Range("N31").Value = -2048.5
$ws.Range("H34").Value = 1156.0845
$ws.Range("I34").Value = 1106.5082
$ws.Range("J34").Value = 1458.5
$ws.Range("K34").Value = 1106.5082
$ws.Range("L34").Value = 1458.5
$ws.Range("M34").Value = -904.5082
$ws.Range("N34").Value = -1862.5
$ws.Range("H106").Value = 29699.5
$ws.Range("J106").Value = 29699.5
$ws.Range("L106").Value = 29699.5
$ws.Range("N106").Value = -32223.5
$ws.Range("H107").Value = 849
$ws.Range("J107").Value = 1637.5
$ws.Range("L107").Value = 1637.5
$ws.Range("N107").Value = -5477.5
$ws.Range("H113").Value = 125001790
$ws.Range("I113").Value = 142858910
$ws.Range("K113").Value = 142858910
$ws.Range("M113").Value = -142856740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 131.3125
$ws.Range("I12").Value = 231.8
$ws.Range("J12").Value = 85.63636
$ws.Range("K12").Value = 695.4000000000001
$ws.Range("L12").Value = 256.90908
$ws.Range("M12").Value = -522.4000000000001
$ws.Range("N12").Value = -602.90908
$ws.Range("H131").Value = 26319308
$ws.Range("I131").Value = 90909500
$ws.Range("J131").Value = 4785
$ws.Range("K131").Value = 272728500
$ws.Range("L131").Value = 14355
$ws.Range("M131").Value = -272723460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 729.9
$ws.Range("I97").Value = 662.375
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 662.375
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -166.375
$ws.Range("N97").Value = -1992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1236.0667
$ws.Range("I93").Value = 1154.6666
$ws.Range("J93").Value = 1561.6666
$ws.Range("K93").Value = 1154.6666
$ws.Range("L93").Value = 1561.6666
$ws.Range("M93").Value = 93.33339999999998
$ws.Range("N93").Value = -4057.6666
$ws.Range("H136").Value = 2017.238
$ws.Range("I136").Value = 1825.2
$ws.Range("J136").Value = 2497.3333
$ws.Range("K136").Value = 5475.6
$ws.Range("L136").Value = 7491.999899999999
$ws.Range("M136").Value = -2925.6
$ws.Range("N136").Value = -12591.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 166670670
$ws.Range("I62").Value = 250001000
$ws.Range("K62").Value = 250001000
$ws.Range("M62").Value = -250000376
$ws.Range("H65").Value = 166670670
$ws.Range("I65").Value = 250001000
$ws.Range("K65").Value = 1250005000
$ws.Range("M65").Value = -1250001880
$ws.Range("H132").Value = 2513.1428
$ws.Range("I132").Value = 2298.6956
$ws.Range("K132").Value = 6896.0868
$ws.Range("M132").Value = -4366.0868
$ws.Range("H136").Value = 1299.125
$ws.Range("I136").Value = 998.26666
$ws.Range("K136").Value = 2994.79998
$ws.Range("M136").Value = -444.7999799999998
